$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row 12
$ws.Range("A12").Value = "eartha51@jolongestr.com"
$ws.Range("B12").Value = "Edward"
$ws.Range("C12").Value = "Artha"
$ws.Range("D12").Value = 6612200748
$ws.Range("E12").Value = "Vendor Admin"

# Row 13
$ws.Range("A13").Value = "pavel@jolongestr.com"
$ws.Range("B13").Value = "Patrick"
$ws.Range("C13").Value = "Avel"
$ws.Range("D13").Value = 6612200748
$ws.Range("E13").Value = "Vendor Report"

# Row 14
$ws.Range("A14").Value = "ethelmae235@jolongestr.com"
$ws.Range("B14").Value = "Edward"
$ws.Range("C14").Value = "Thelmae"
$ws.Range("D14").Value = 6612200748
$ws.Range("E14").Value = "Vendor Admin"

$ws.Hyperlinks.Add($ws.Range("A12"), "mailto:eartha51@jolongestr.com")
$ws.Hyperlinks.Add($ws.Range("A13"), "mailto:pavel@jolongestr.com")
$ws.Hyperlinks.Add($ws.Range("A14"), "mailto:ethelmae235@jolongestr.com")

$ws.Range("A12").Style = "Hyperlink"
$ws.Range("A13").Style = "Hyperlink"
$ws.Range("A14").Style = "Hyperlink"
